$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A82").Value = "'2025/10/09"
$ws.Range("A82").Style = "Normal"
$ws.Range("B82").Value = "木"
$ws.Range("C82").Value = 7
$ws.Range("D82").Value = 16
